# Add a new row 6 that duplicates row 3's submission (same respondent
# resubmitting), right below the existing row 5, mirroring the row-5
# formatting. This matches an Excel "copy row -> insert copied cells"
# user action (copy row 3, insert the copy at row 6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(3).Copy()
$ws.Rows.Item(6).Insert()

# Re-enter the submission timestamp precisely (guards against any float
# rounding introduced by the copy/insert operation).
$ws.Range("C6").Value = 45426.769814814812

# Give the new row the same explicit "custom" row height the other
# data rows carry (15.75pt).
$ws.Rows.Item(6).RowHeight = 15.75

# The header/data rows above settle to their natural (non-custom) height
# once the sheet is re-laid out after the insert.
$ws.Rows.Item(1).RowHeight = 12.75
$ws.Rows.Item(2).RowHeight = 12.75
$ws.Rows.Item(3).RowHeight = 12.75
$ws.Rows.Item(4).RowHeight = 12.75

# Move the selection to the newly added row, like Excel leaves the
# pasted-into range selected.
$ws.Range("A6:V6").Select() | Out-Null
